# "Final changes, with exception handling"
# Update the Doctors sheet: rows 2-6 (Dentist entries from Delhi) are replaced
# with new Dermatologist entries from Mumbai. Only columns A (Doctors Names),
# B (Field), C (Experience) and D (Practise Location) change; column E
# (Surgeries List) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

$data = @(
    @("Dr. Kiran Katkar", "Dermatologist", "39 years experience overall", "Dadar West,Mumbai"),
    @("Dr. Sonia Tekchandani", "Dermatologist", "38 years experience overall", "Andheri West,Mumbai"),
    @("Dr. Jolly Shah Kapadia", "Dermatologist", "25 years experience overall", "Mulund West,Mumbai"),
    @("Dr. Resham Vasani Bhojani", "Dermatologist", "25 years experience overall", "Matunga,Mumbai"),
    @("Dr. Zeenat Bhalwani", "Dermatologist", "23 years experience overall", "Andheri West,Mumbai")
)

try {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $row = $i + 2
        $values = $data[$i]
        $ws.Cells.Item($row, 1).Value = $values[0]
        $ws.Cells.Item($row, 2).Value = $values[1]
        $ws.Cells.Item($row, 3).Value = $values[2]
        $ws.Cells.Item($row, 4).Value = $values[3]
    }
}
catch {
    Write-Host "Error while updating Doctors sheet: $_"
}
